$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SNMP")

$ws.Range("F1").Value = "SNMPserver"
$ws.Range("G1").Value = "SNMPuser"
$ws.Range("H1").Value = "SNMPpass"

$ws.Range("B2").Value = "MD5|SHA"
$ws.Range("F2").Value = "10.30.4.77"
$ws.Range("H2").Value = "rootpw"
$ws.Range("D2").Value = "DES|AES128"

$ws.Range("E9").Select()
